$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Update category labels (column A)
$ws.Range("A3").Value = "Kun fuglebeskyt"
$ws.Range("A4").Value = "Kun habitatomraade"
$ws.Range("A5").Value = "Habitatatomr. og fuglebesk."
$ws.Range("A6").Value = "Habitatnaturtype"
$ws.Range("A7").Value = "Ramsar"
$ws.Range("A8").Value = "Havstrategi standard"
$ws.Range("A9").Value = "Havstrategi streng"
$ws.Range("A10").Value = "Vildt Reservater"

# Update Area_Sq_Km (column B) and Proportion (column C) values
$ws.Range("B2").Value = 26284.8191
$ws.Range("C2").Value = 25.0282495473987

$ws.Range("B3").Value = 22270.9031
$ws.Range("C3").Value = 21.2062224325042

$ws.Range("B4").Value = 16074.5452
$ws.Range("C4").Value = 15.3060870267332

$ws.Range("B5").Value = 12060.63
$ws.Range("C5").Value = 11.484060673594

$ws.Range("B6").Value = 7105.5811
$ws.Range("C6").Value = 6.76589236827123

$ws.Range("B7").Value = 5647.0395
$ws.Range("C7").Value = 5.37707766876043

$ws.Range("B8").Value = 3494.8231
$ws.Range("C8").Value = 3.32774992051642

$ws.Range("B9").Value = 3280.8196
$ws.Range("C9").Value = 3.12397705140747

$ws.Range("B10").Value = 2290.1955
$ws.Range("C10").Value = 2.18071063256165
